# Auto-generated: apply 2023-09-17 violent crime YTD data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 5441
$ws.Cells.Item(3, 10).Value = 5789
$ws.Cells.Item(4, 10).Value = 1264
$ws.Cells.Item(6, 10).Value = 7230
$ws.Cells.Item(7, 10).Value = 20171
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(4, 10).Value = 86
$ws.Cells.Item(7, 10).Value = 583
$ws.Cells.Item(8, 10).Value = 1266
$ws.Cells.Item(10, 10).Value = 136
$ws.Cells.Item(11, 10).Value = 315
$ws.Cells.Item(15, 10).Value = 222
$ws.Cells.Item(18, 10).Value = 171
$ws.Cells.Item(19, 10).Value = 584
$ws.Cells.Item(22, 10).Value = 53
$ws.Cells.Item(27, 10).Value = 122
$ws.Cells.Item(29, 10).Value = 1130
$ws.Cells.Item(31, 10).Value = 183
$ws.Cells.Item(33, 10).Value = 930
$ws.Cells.Item(34, 10).Value = 98
$ws.Cells.Item(36, 10).Value = 276
$ws.Cells.Item(37, 10).Value = 622
$ws.Cells.Item(41, 10).Value = 126
$ws.Cells.Item(42, 10).Value = 832
$ws.Cells.Item(44, 10).Value = 153
$ws.Cells.Item(48, 10).Value = 232
$ws.Cells.Item(49, 10).Value = 135
$ws.Cells.Item(50, 10).Value = 127
$ws.Cells.Item(51, 10).Value = 250
$ws.Cells.Item(52, 10).Value = 516
$ws.Cells.Item(53, 10).Value = 274
$ws.Cells.Item(54, 10).Value = 386
$ws.Cells.Item(55, 10).Value = 267
$ws.Cells.Item(63, 10).Value = 68
$ws.Cells.Item(65, 10).Value = 518
$ws.Cells.Item(67, 10).Value = 766
$ws.Cells.Item(69, 10).Value = 48
$ws.Cells.Item(70, 10).Value = 28
$ws.Cells.Item(72, 10).Value = 79
$ws.Cells.Item(75, 10).Value = 61
$ws.Cells.Item(76, 10).Value = 291
$ws.Cells.Item(79, 10).Value = 578
$ws.Cells.Item(84, 10).Value = 173
$ws.Cells.Item(85, 10).Value = 857
$ws.Cells.Item(88, 10).Value = 220
$ws.Cells.Item(89, 10).Value = 262
$ws.Cells.Item(92, 10).Value = 61
$ws.Cells.Item(93, 10).Value = 91
$ws.Cells.Item(94, 10).Value = 205
$ws.Cells.Item(96, 10).Value = 241
$ws.Cells.Item(97, 10).Value = 166
$ws.Cells.Item(101, 10).Value = 20171
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(2, 10).Value = 71
$ws.Cells.Item(4, 10).Value = 15
$ws.Cells.Item(7, 10).Value = 241
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 183
$ws.Cells.Item(7, 10).Value = 583
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 10).Value = 97
$ws.Cells.Item(3, 10).Value = 66
$ws.Cells.Item(6, 10).Value = 126
$ws.Cells.Item(7, 10).Value = 315
$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(4, 10).Value = 28
$ws.Cells.Item(6, 10).Value = 75
$ws.Cells.Item(7, 10).Value = 262
$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 10).Value = 223
$ws.Cells.Item(6, 10).Value = 247
$ws.Cells.Item(7, 10).Value = 857
$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 10).Value = 123
$ws.Cells.Item(3, 10).Value = 158
$ws.Cells.Item(6, 10).Value = 206
$ws.Cells.Item(7, 10).Value = 516
$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Cells.Item(4, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 48
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 10).Value = 54
$ws.Cells.Item(6, 10).Value = 173
$ws.Cells.Item(7, 10).Value = 274
$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 10).Value = 350
$ws.Cells.Item(3, 10).Value = 385
$ws.Cells.Item(4, 10).Value = 72
$ws.Cells.Item(6, 10).Value = 424
$ws.Cells.Item(7, 10).Value = 1266
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 10).Value = 227
$ws.Cells.Item(3, 10).Value = 304
$ws.Cells.Item(4, 10).Value = 39
$ws.Cells.Item(6, 10).Value = 320
$ws.Cells.Item(7, 10).Value = 930
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 10).Value = 214
$ws.Cells.Item(6, 10).Value = 178
$ws.Cells.Item(7, 10).Value = 622
$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(3, 10).Value = 147
$ws.Cells.Item(7, 10).Value = 518
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 10).Value = 70
$ws.Cells.Item(6, 10).Value = 49
$ws.Cells.Item(7, 10).Value = 183
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 10).Value = 294
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(7, 10).Value = 766
$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(3, 10).Value = 54
$ws.Cells.Item(7, 10).Value = 173
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(4, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 135
$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(3, 10).Value = 74
$ws.Cells.Item(7, 10).Value = 386
$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 10).Value = 337
$ws.Cells.Item(3, 10).Value = 389
$ws.Cells.Item(4, 10).Value = 63
$ws.Cells.Item(6, 10).Value = 298
$ws.Cells.Item(7, 10).Value = 1130
$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 10).Value = 36
$ws.Cells.Item(3, 10).Value = 43
$ws.Cells.Item(7, 10).Value = 232
$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 10).Value = 147
$ws.Cells.Item(6, 10).Value = 213
$ws.Cells.Item(7, 10).Value = 584
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(6, 10).Value = 61
$ws.Cells.Item(7, 10).Value = 153
$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(6, 10).Value = 161
$ws.Cells.Item(7, 10).Value = 291
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(2, 10).Value = 29
$ws.Cells.Item(7, 10).Value = 126
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(6, 10).Value = 425
$ws.Cells.Item(7, 10).Value = 832
$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(3, 10).Value = 27
$ws.Cells.Item(7, 10).Value = 136
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 10).Value = 64
$ws.Cells.Item(4, 10).Value = 10
$ws.Cells.Item(6, 10).Value = 129
$ws.Cells.Item(7, 10).Value = 267
$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 10).Value = 163
$ws.Cells.Item(3, 10).Value = 205
$ws.Cells.Item(6, 10).Value = 160
$ws.Cells.Item(7, 10).Value = 578
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(3, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 171
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(3, 10).Value = 90
$ws.Cells.Item(6, 10).Value = 80
$ws.Cells.Item(7, 10).Value = 276
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(3, 10).Value = 29
$ws.Cells.Item(6, 10).Value = 33
$ws.Cells.Item(7, 10).Value = 91
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 98
$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(2, 10).Value = 36
$ws.Cells.Item(3, 10).Value = 41
$ws.Cells.Item(6, 10).Value = 112
$ws.Cells.Item(7, 10).Value = 205
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 10).Value = 64
$ws.Cells.Item(7, 10).Value = 222
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(4, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 127
$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(6, 10).Value = 115
$ws.Cells.Item(7, 10).Value = 166
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(3, 10).Value = 23
$ws.Cells.Item(7, 10).Value = 61
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(2, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 28
$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(6, 10).Value = 101
$ws.Cells.Item(7, 10).Value = 220
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(3, 10).Value = 30
$ws.Cells.Item(7, 10).Value = 122
$ws = $wb.Worksheets.Item('Pullman')
$ws.Cells.Item(3, 10).Value = 17
$ws.Cells.Item(7, 10).Value = 61
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 10).Value = 59
$ws.Cells.Item(7, 10).Value = 250
$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(6, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 53
$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(6, 10).Value = 28
$ws.Cells.Item(7, 10).Value = 79
$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(2, 10).Value = 29
$ws.Cells.Item(3, 10).Value = 20
$ws.Cells.Item(6, 10).Value = 31
$ws.Cells.Item(7, 10).Value = 86

